# Trade #7 closed at 2026-02-16 22:57:57 - base_strategy DOWN +0.000%
#
# The bot logged a new (still-open) trade as row 8 on both the "All Trades"
# log sheet and the per-strategy "base_strategy" sheet. Append the identical
# row to each.

$wb = $excel.ActiveWorkbook

$targetRow = 8

$tradeNum    = 7
$tradeDate   = "2026-02-16"
$tradeTime   = "22:57:57"
$strategy    = "base_strategy"
$side        = "DOWN"
$entryPrice  = 0.5
$exitPrice   = ""
$status      = "OPEN"
$pnlPct      = 0
$pnlDollar   = 0
$capAfter    = 100
$entrySlip   = 0
$exitSlip    = 0
$confidence  = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason  = ""
$durationMin = 0

foreach ($sheetName in @("All Trades", "base_strategy")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($targetRow, 1).Value = $tradeNum

    # Write the date as literal text ("2026-02-16"), not an auto-converted
    # date serial: format the cell as Text first, then reset the style back
    # to Normal afterwards so no stray number format sticks to the cell.
    $ws.Cells.Item($targetRow, 2).NumberFormat = "@"
    $ws.Cells.Item($targetRow, 2).Value = $tradeDate
    $ws.Cells.Item($targetRow, 2).Style = "Normal"

    $ws.Cells.Item($targetRow, 3).Value = $tradeTime
    $ws.Cells.Item($targetRow, 4).Value = $strategy
    $ws.Cells.Item($targetRow, 5).Value = $side
    $ws.Cells.Item($targetRow, 6).Value = $entryPrice
    $ws.Cells.Item($targetRow, 7).Value = $exitPrice
    $ws.Cells.Item($targetRow, 8).Value = $status
    $ws.Cells.Item($targetRow, 9).Value = $pnlPct
    $ws.Cells.Item($targetRow, 10).Value = $pnlDollar
    $ws.Cells.Item($targetRow, 11).Value = $capAfter
    $ws.Cells.Item($targetRow, 12).Value = $entrySlip
    $ws.Cells.Item($targetRow, 13).Value = $exitSlip
    $ws.Cells.Item($targetRow, 14).Value = $confidence
    $ws.Cells.Item($targetRow, 15).Value = $entryReason
    $ws.Cells.Item($targetRow, 16).Value = $exitReason
    $ws.Cells.Item($targetRow, 17).Value = $durationMin
}
